# Apply the "#5: insurance, claim, debt, investment done" edit.
# This brings the 債務 (debt) and 事業投資 (business investment) sheets in
# line with the other sheets' 14-column schema: proper header row (B..N)
# plus the common trailing columns (property_category, category, date,
# legislator_name, legislator_id, source_file, index) on every data row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "債務" (debt)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("債務")

# Header row (row 1): was a duplicate of row 2's data, now real column names.
$ws.Range("B1").Value = "species"
$ws.Range("C1").Value = "debtor"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "total"
$ws.Range("F1").Value = "register_date"
$ws.Range("G1").Value = "register_reason"
$ws.Range("H1").Value = "property_category"
$ws.Range("I1").Value = "category"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "legislator_name"
$ws.Range("L1").Value = "legislator_id"
$ws.Range("M1").Value = "source_file"
$ws.Range("N1").Value = "index"

# Row 2 (index 129)
$ws.Range("B2").Value = "房屋貸款"
$ws.Range("C2").Value = "葉宜津"
$ws.Range("D2").Value = "華南商業銀行臺南市新營區新進路2段109號"
$ws.Range("E2").Value = 40000000
$ws.Range("F2").Value = "102年06月03日"
$ws.Range("G2").Value = "借貸"
$ws.Range("H2").Value = "debt"
$ws.Range("I2").Value = "normal"
$ws.Range("J2").Value = "2013-11-08"
$ws.Range("K2").Value = "葉宜津"
$ws.Range("L2").Value = 855
$ws.Range("M2").Value = "tmpabd41"
$ws.Range("N2").Value = 129

# Row 3 (index 130)
$ws.Range("B3").Value = "長期擔保貸款"
$ws.Range("C3").Value = "趙哲宏"
$ws.Range("D3").Value = "華南銀行新營分行臺南市新營區新進路2段109號"
$ws.Range("F3").Value = "100年04月15日"
$ws.Range("G3").Value = "借貸"
$ws.Range("H3").Value = "debt"
$ws.Range("I3").Value = "normal"
$ws.Range("J3").Value = "2013-11-08"
$ws.Range("K3").Value = "葉宜津"
$ws.Range("L3").Value = 855
$ws.Range("M3").Value = "tmpabd41"
$ws.Range("N3").Value = 130

# Row 4 (index 131)
$ws.Range("B4").Value = "私人債務"
$ws.Range("C4").Value = "葉宜津"
$ws.Range("D4").Value = "葉光彰臺北市松山區民生東路"
$ws.Range("F4").Value = "102年06月03日"
$ws.Range("G4").Value = "借貸"
$ws.Range("H4").Value = "debt"
$ws.Range("I4").Value = "normal"
$ws.Range("J4").Value = "2013-11-08"
$ws.Range("K4").Value = "葉宜津"
$ws.Range("L4").Value = 855
$ws.Range("M4").Value = "tmpabd41"
$ws.Range("N4").Value = 131

# ---------------------------------------------------------------
# Sheet "事業投資" (business investment)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("事業投資")

# Header row (row 1): was a duplicate of row 2's data, now real column names.
$ws2.Range("B1").Value = "owner"
$ws2.Range("C1").Value = "company"
$ws2.Range("D1").Value = "address"
$ws2.Range("E1").Value = "total"
$ws2.Range("F1").Value = "register_date"
$ws2.Range("G1").Value = "register_reason"
$ws2.Range("H1").Value = "property_category"
$ws2.Range("I1").Value = "category"
$ws2.Range("J1").Value = "date"
$ws2.Range("K1").Value = "legislator_name"
$ws2.Range("L1").Value = "legislator_id"
$ws2.Range("M1").Value = "source_file"
$ws2.Range("N1").Value = "index"

# Row 2 (index 136)
$ws2.Range("H2").Value = "investment"
$ws2.Range("I2").Value = "normal"
$ws2.Range("J2").Value = "2013-11-08"
$ws2.Range("K2").Value = "葉宜津"
$ws2.Range("L2").Value = 855
$ws2.Range("M2").Value = "tmpabd41"
$ws2.Range("N2").Value = 136
